# Applies the "3 routes added" change to the queries sheet:
#  - Fixes G18 (weight_tons) so it is stored as a real number (12) instead of text.
#  - Appends 6 new quote rows (19-24) with data from the new routes.
#  - Row 24's weight_tons (G24) stays a text value of "12" (matches source diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing row 18: weight_tons should be numeric, not text ---
$ws.Cells.Item(18, 7).Value = 12

# --- New row data (rows 19-24) ---
# columns: A quote_id, B company_name, C shipping_from, D destination,
#          E commodity, F cargo_type, G weight_tons, H num_containers,
#          I timestamp, J container_type
$newRows = @(
    @("QUOTE-20251223073948", "nabeel", "Karachi Port", "Almaty",       "Food Item", "", 12, "", "2025-12-23 07:39:48", "Dry"),
    @("QUOTE-20251223074946", "nabeel", "Almaty",       "Karachi Port", "Food Item", "", 21, "", "2025-12-23 07:49:46", "Dry"),
    @("QUOTE-20251223080055", "nabeel", "Almaty",       "Karachi Port", "Food Item", "", 12, "", "2025-12-23 08:00:55", "Dry"),
    @("QUOTE-20251223080130", "nabeel", "Almaty",       "Karachi Port", "Food Item", "", 12, "", "2025-12-23 08:01:30", "Dry"),
    @("QUOTE-20251223082457", "nabeel", "Karachi Port", "almaty",       "Food Item", "", 12, "", "2025-12-23 08:24:57", "Dry"),
    @("QUOTE-20251223082807", "nabeel", "Karachi Port", "almaty",       "Food Item", "", 12, "", "2025-12-23 08:28:07", "Dry")
)

$startRow = 19
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 7).Value = $data[6]
    $ws.Cells.Item($r, 9).Value = $data[8]
    $ws.Cells.Item($r, 10).Value = $data[9]
}

# Row 24's weight_tons stays textual "12" (not converted to a number),
# matching the source data for that route. Force text by prefixing with an
# apostrophe, then reset the style picked up from the quote-prefix so no
# stray formatting is introduced.
$ws.Cells.Item(24, 7).Value = "'12"
$ws.Cells.Item(24, 7).Style = "Normal"
